# Applies the "INPC subyacente/no subyacente" disaggregation edit:
# Inserts 4 new rows (2 under "Subyacente", 2 under "No subyacente")
# into the series table on Sheet1, pushing the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the two new rows after "Subyacente" (old row 4, before "INPC CCM") ---
$ws.Range("A4:A5").EntireRow.Insert()

# --- Insert the two new rows after "No subyacente" (now at row 7, before "INPC quincenal") ---
$ws.Range("A8:A9").EntireRow.Insert()

# The "idEstructura" column (B) holds digit-only strings that must stay
# text (matching the rest of the column), so force text formatting before
# assigning, then restore the default style to match sibling cells.
$idEstructuraCells = @("B4", "B5", "B8", "B9")
foreach ($addr in $idEstructuraCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Fill in the new row 4: Subyacente - Mercancias ---
$ws.Range("A4").Value = "Subyacente - Mercancias"
$ws.Range("B4").Value = "112001700010"
$ws.Range("C4").Value = "e|865548"
$ws.Range("D4").Value = "v_subyacente_mercancias"

# --- Fill in the new row 5: Subyacente - Servicios ---
$ws.Range("A5").Value = "Subyacente - Servicios"
$ws.Range("B5").Value = "112001700010"
$ws.Range("C5").Value = "e|865551"
$ws.Range("D5").Value = "v_subyacente_servicios"

# --- Fill in the new row 8: No subyacente - Agropecuarios ---
$ws.Range("A8").Value = "No subyacente - Agropecuarios"
$ws.Range("B8").Value = "112001700010"
$ws.Range("C8").Value = "e|865556"
$ws.Range("D8").Value = "v_nsubyacente_agropecuarios"

# --- Fill in the new row 9: No subyacente - Energéticos y tarifas autorizadas ---
$ws.Range("A9").Value = "No subyacente - Energéticos y tarifas autorizadas"
$ws.Range("B9").Value = "112001700010"
$ws.Range("C9").Value = "e|865559"
$ws.Range("D9").Value = "v_nsubyacente_energeticos"

# Restore the default (unstyled) cell style on the idEstructura cells so
# they match their sibling data cells, which carry no explicit style.
foreach ($addr in $idEstructuraCells) {
    $ws.Range($addr).Style = "Normal"
}
